$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.768.25"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.298.45"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.20%  "
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.55%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.05%  "
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "2.657.29"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "2.311.16"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.796"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "42.716.95"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0695"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").Value = "1.961.76"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "2.523.60"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("E51").Value = "  -2.82%  "
